# Update simulation results in Sheet1 (pl_mw.xlsx) for the 380 kV case.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.096718192662536
$ws.Range("C2").Value = 0.1042318063352177
$ws.Range("D2").Value = 0.5500440473496155
$ws.Range("E2").Value = 0.1859857314196667
$ws.Range("G2").Value = 0.002576913897089775
$ws.Range("J2").Value = 0.07349989187674222
$ws.Range("K2").Value = 0.551167743408314
$ws.Range("L2").Value = 0.4023042471223448
$ws.Range("O2").Value = 8.37175143879864
$ws.Range("B3").Value = 1.064984540794825
$ws.Range("C3").Value = 0.1034372542101671
$ws.Range("D3").Value = 0.5475900088491557
$ws.Range("E3").Value = 0.1860915640303418
$ws.Range("G3").Value = 0.002580121558918526
$ws.Range("J3").Value = 0.07380485733061981
$ws.Range("K3").Value = 0.5222684740820682
$ws.Range("L3").Value = 0.3988265014650807
$ws.Range("O3").Value = 8.392468579742797
$ws.Range("B4").Value = 1.046023487550912
$ws.Range("C4").Value = 0.1029404290832119
$ws.Range("D4").Value = 0.5463102605992987
$ws.Range("E4").Value = 0.1862201825593388
$ws.Range("G4").Value = 0.002582197010429348
$ws.Range("J4").Value = 0.07400550936675998
$ws.Range("K4").Value = 0.5047806029703281
$ws.Range("L4").Value = 0.3968622167576825
$ws.Range("O4").Value = 8.408341354318026
$ws.Range("B5").Value = 1.038428955501075
$ws.Range("C5").Value = 0.1027357212078925
$ws.Range("D5").Value = 0.5458459537376541
$ws.Range("E5").Value = 0.1862886283748288
$ws.Range("G5").Value = 0.002583069492951905
$ws.Range("J5").Value = 0.07409065352559452
$ws.Range("K5").Value = 0.4977190565096095
$ws.Range("L5").Value = 0.3961048805959422
$ws.Range("O5").Value = 8.415602730571663
$ws.Range("B6").Value = 1.037175894720519
$ws.Range("C6").Value = 0.1027015940851363
$ws.Range("D6").Value = 0.5457723139672055
$ws.Range("E6").Value = 0.1863009629924939
$ws.Range("G6").Value = 0.002583215984058089
$ws.Range("J6").Value = 0.07410499579948748
$ws.Range("K6").Value = 0.4965504248062587
$ws.Range("L6").Value = 0.3959817334479752
$ws.Range("O6").Value = 8.416856389326256
$ws.Range("B7").Value = 1.045920528753186
$ws.Range("C7").Value = 0.1029376774128998
$ws.Range("D7").Value = 0.5463037670438666
$ws.Range("E7").Value = 0.1862210406866289
$ws.Range("G7").Value = 0.002582208668703925
$ws.Range("J7").Value = 0.07400664396928125
$ws.Range("K7").Value = 0.5046851051228316
$ws.Range("L7").Value = 0.3968518282955742
$ws.Range("O7").Value = 8.408436072063637
$ws.Range("B8").Value = 1.085668117724623
$ws.Range("C8").Value = 0.1039597099179481
$ws.Range("D8").Value = 0.5491508361116502
$ws.Range("E8").Value = 0.1860090327510058
$ws.Range("G8").Value = 0.00257799796113111
$ws.Range("J8").Value = 0.07360226645456081
$ws.Range("K8").Value = 0.5411503069018408
$ws.Range("L8").Value = 0.4010696797347606
$ws.Range("O8").Value = 8.378240666267885
$ws.Range("B9").Value = 1.167745812952973
$ws.Range("C9").Value = 0.1058925313429881
$ws.Range("D9").Value = 0.5565314412555722
$ws.Range("E9").Value = 0.1860970072762065
$ws.Range("G9").Value = 0.002570577614984969
$ws.Range("J9").Value = 0.07291531913700222
$ws.Range("K9").Value = 0.6146782439615208
$ws.Range("L9").Value = 0.4106942163902971
$ws.Range("O9").Value = 8.344027698314363
$ws.Range("B10").Value = 1.230546988441375
$ws.Range("C10").Value = 0.1072688693615547
$ws.Range("D10").Value = 0.5630457229647874
$ws.Range("E10").Value = 0.1864673073681651
$ws.Range("G10").Value = 0.002565630860291538
$ws.Range("J10").Value = 0.07247484559618478
$ws.Range("K10").Value = 0.6699172145327736
$ws.Range("L10").Value = 0.4185863049101783
$ws.Range("O10").Value = 8.334123209951429
$ws.Range("B11").Value = 1.259655574632177
$ws.Range("C11").Value = 0.1078854803065923
$ws.Range("D11").Value = 0.5662455299896862
$ws.Range("E11").Value = 0.1867018544887742
$ws.Range("G11").Value = 0.002563489000323902
$ws.Range("J11").Value = 0.07228832122894424
$ws.Range("K11").Value = 0.695308941200949
$ws.Range("L11").Value = 0.4223541023460626
$ws.Range("O11").Value = 8.332923581390901
$ws.Range("B12").Value = 1.270755408177962
$ws.Range("C12").Value = 0.1081176047726728
$ws.Range("D12").Value = 0.567491122828585
$ws.Range("E12").Value = 0.1868001488048883
$ws.Range("G12").Value = 0.002562693443914557
$ws.Range("J12").Value = 0.07221967404313112
$ws.Range("K12").Value = 0.7049616476256517
$ws.Range("L12").Value = 0.4238063249238877
$ws.Range("O12").Value = 8.332944486669476
$ws.Range("B13").Value = 1.268361443557751
$ws.Range("C13").Value = 0.1080676737889377
$ws.Range("D13").Value = 0.5672213562480692
$ws.Range("E13").Value = 0.1867785583450257
$ws.Range("G13").Value = 0.002562864092258245
$ws.Range("J13").Value = 0.07223437022963708
$ws.Range("K13").Value = 0.7028811077667569
$ws.Range("L13").Value = 0.4234924330618384
$ws.Range("O13").Value = 8.332918854054014
$ws.Range("B14").Value = 1.260567224029302
$ws.Range("C14").Value = 0.1079046048815968
$ws.Range("D14").Value = 0.5663473269342489
$ws.Range("E14").Value = 0.1867097514602243
$ws.Range("G14").Value = 0.002563423238621975
$ws.Range("J14").Value = 0.07228263382470956
$ws.Range("K14").Value = 0.6961023279320102
$ws.Range("L14").Value = 0.4224730683311577
$ws.Range("O14").Value = 8.332915779223015
$ws.Range("B15").Value = 1.255803054835553
$ws.Range("C15").Value = 0.1078045414710189
$ws.Range("D15").Value = 0.5658163698801388
$ws.Range("E15").Value = 0.1866688385340609
$ws.Range("G15").Value = 0.002563767750988721
$ws.Range("J15").Value = 0.07231245507005468
$ws.Range("K15").Value = 0.6919549884117657
$ws.Range("L15").Value = 0.421851987835268
$ws.Range("O15").Value = 8.332975771709016
$ws.Range("B16").Value = 1.228655488159802
$ws.Range("C16").Value = 0.1072283807129395
$ws.Range("D16").Value = 0.5628413567957011
$ws.Range("E16").Value = 0.186453306774812
$ws.Range("G16").Value = 0.002565773011723266
$ws.Range("J16").Value = 0.07248731357415927
$ws.Range("K16").Value = 0.6682630646191114
$ws.Range("L16").Value = 0.4183436361258543
$ws.Range("O16").Value = 8.334268114562263
$ws.Range("B17").Value = 1.212139211285205
$ws.Range("C17").Value = 0.1068724880156466
$ws.Range("D17").Value = 0.5610767642619692
$ws.Range("E17").Value = 0.1863379937105343
$ws.Range("G17").Value = 0.002567030896040118
$ws.Range("J17").Value = 0.07259812647400388
$ws.Range("K17").Value = 0.6537959429718683
$ws.Range("L17").Value = 0.4162368000341843
$ws.Range("O17").Value = 8.335907509148626
$ws.Range("B18").Value = 1.202690350806222
$ws.Range("C18").Value = 0.1066668949598863
$ws.Range("D18").Value = 0.560084077107291
$ws.Range("E18").Value = 0.1862778916017191
$ws.Range("G18").Value = 0.002567764609414571
$ws.Range("J18").Value = 0.07266316703429787
$ws.Range("K18").Value = 0.6454996411200113
$ws.Range("L18").Value = 0.4150417308927814
$ws.Range("O18").Value = 8.337161632029421
$ws.Range("B19").Value = 1.199499883786302
$ws.Range("C19").Value = 0.1065971315501031
$ws.Range("D19").Value = 0.5597517965689605
$ws.Range("E19").Value = 0.1862586117097322
$ws.Range("G19").Value = 0.00256801478834956
$ws.Range("J19").Value = 0.0726854127808334
$ws.Range("K19").Value = 0.6426949295090481
$ws.Range("L19").Value = 0.4146399771743461
$ws.Range("O19").Value = 8.337639709735157
$ws.Range("B20").Value = 1.213892137234126
$ws.Range("C20").Value = 0.1069104658799063
$ws.Range("D20").Value = 0.5612623050890164
$ws.Range("E20").Value = 0.1863496251386216
$ws.Range("G20").Value = 0.002566895935714852
$ws.Range("J20").Value = 0.07258619534518296
$ws.Range("K20").Value = 0.6553334288812778
$ws.Range("L20").Value = 0.4164593458946371
$ws.Range("O20").Value = 8.33570078771956
$ws.Range("B21").Value = 1.262854489139045
$ws.Range("C21").Value = 0.1079525394890268
$ws.Range("D21").Value = 0.5666031314371196
$ws.Range("E21").Value = 0.1867297047463587
$ws.Range("G21").Value = 0.002563258583325582
$ws.Range("J21").Value = 0.07226840379729715
$ws.Range("K21").Value = 0.6980924080937143
$ws.Range("L21").Value = 0.4227717908841413
$ws.Range("O21").Value = 8.332903787872226
$ws.Range("B22").Value = 1.295303005644996
$ws.Range("C22").Value = 0.1086255900991091
$ws.Range("D22").Value = 0.5702911909534123
$ws.Range("E22").Value = 0.1870333315861288
$ws.Range("G22").Value = 0.002560971790286709
$ws.Range("J22").Value = 0.0720722791825299
$ws.Range("K22").Value = 0.7262557269131378
$ws.Range("L22").Value = 0.4270455704949256
$ws.Range("O22").Value = 8.333845393108675
$ws.Range("B23").Value = 1.277943758076674
$ws.Range("C23").Value = 0.1082671052774984
$ws.Range("D23").Value = 0.5683047648309127
$ws.Range("E23").Value = 0.1868662365760017
$ws.Range("G23").Value = 0.002562184044848863
$ws.Range("J23").Value = 0.07217589783477507
$ws.Range("K23").Value = 0.7112046456341545
$ws.Range("L23").Value = 0.4247510455042089
$ws.Range("O23").Value = 8.333089492177237
$ws.Range("B24").Value = 1.213099494314974
$ws.Range("C24").Value = 0.1068932991567593
$ws.Range("D24").Value = 0.5611783541682058
$ws.Range("E24").Value = 0.1863443472778279
$ws.Range("G24").Value = 0.00256695691848015
$ws.Range("J24").Value = 0.07259158525678977
$ws.Range("K24").Value = 0.6546382659343806
$ws.Range("L24").Value = 0.4163586825086441
$ws.Range("O24").Value = 8.335793275772602
$ws.Range("B25").Value = 1.145101339317961
$ws.Range("C25").Value = 0.105377314046649
$ws.Range("D25").Value = 0.5543427773223613
$ws.Range("E25").Value = 0.1860194279736405
$ws.Range("G25").Value = 0.002572495964101908
$ws.Range("J25").Value = 0.07308984784794426
$ws.Range("K25").Value = 0.5945721458484741
$ws.Range("L25").Value = 0.4079460649279464
$ws.Range("O25").Value = 8.350607717548939
